# Daily attendance processing - normalise "Recorded By" (column G) ordering.
# For every row whose "Recorded By" text is a comma-separated list that
# starts with a "System" token (any case), move that leading token to the
# end of the list (re-emitted as "System") so the real recorder appears
# first. Single-valued cells, and cells that don't start with "System",
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Length -lt 2) {
        continue
    }

    if ($parts[0].Trim().ToLower() -ne "system") {
        continue
    }

    $lastIdx = $parts.Length - 1
    $first = $parts[0]
    $last = $parts[$lastIdx]

    $parts[0] = $last
    $parts[$lastIdx] = "System"

    $newText = $parts -join ", "
    $cell.Value = $newText
}
